$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting down for the three new rows (copy the previous data row's
# look: A/B plain style, C "Neutral" style, D "Neutral" style)
$ws.Range("A49:D49").Copy()
$ws.Range("A50:D52").PasteSpecial(-4122)

# Row 51's "My Status" cell (D) should use the "Good" style instead, matching
# row 2's D cell
$ws.Range("D2").Copy()
$ws.Range("D51").PasteSpecial(-4122)

# Fill in the new problems, in the order they were typed so the shared
# string table grows in the same order as the source workbook
$ws.Cells.Item(51, 2).Value = "Letter Combinations of a Phone Number"
$ws.Cells.Item(52, 2).Value = "4Sum"
$ws.Cells.Item(50, 2).Value = "3Sum Closest"

$ws.Cells.Item(50, 1).Value = 16
$ws.Cells.Item(51, 1).Value = 17
$ws.Cells.Item(52, 1).Value = 18

# Reflect the user's final scroll position / selection after adding the rows
$ws.Range("F53").Select() | Out-Null
